$d = $word.ActiveDocument

# Walk forward through the document with sequential Find calls (each call searches onward from
# where the previous match ended) so we land on the exact Range boundaries of the runs involved,
# instead of hard-coding character offsets.

$r = $d.Content
$null = $r.Find.Execute("poser ", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

$null = $r.Find.Execute("<corr>", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$corrStart = $r.Start

$null = $r.Find.Execute(" sur", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$surStart = $r.Start

# Grab just the "s" character (the 2nd character of the found " sur", i.e. surStart+1..surStart+2)
# together with its run's own formatting (that run carries only rtl=0 - no font/color/size
# override), so splicing it in elsewhere keeps that exact formatting instead of inheriting from
# its new neighbours.
$sChar = $d.Range($surStart + 1, $surStart + 2)

# Splice a copy of that formatted "s" in right before "<corr>".
$insertionPoint = $d.Range($corrStart, $corrStart)
$insertionPoint.FormattedText = $sChar.FormattedText

# The previously found " sur" run has shifted right by one character because of the insert above.
# Remove its leading " s" (now at surStart+1 .. surStart+3) so it reads "ur".
$toDelete = $d.Range($surStart + 1, $surStart + 3)
$toDelete.Text = ""
